# Apply updates to the CLX balance sheet data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CLX")

# Row 4: Inventory
$ws.Range("C4").Value = 609000000.0
$ws.Range("D4").Value = 534000000.0
$ws.Range("F4").Value = 457000000.0
$ws.Range("G4").Value = 514000000.0

# Row 14: Accounts Payable
$ws.Range("C14").Value = 1373000000.0
$ws.Range("D14").Value = 1391000000.0
$ws.Range("E14").Value = 1322000000.0
$ws.Range("F14").Value = 1083000000.0
$ws.Range("G14").Value = 942000000.0

# Row 21: Long Term Tax Liability (Deferred)
$ws.Range("C21").Value = 128000000.0
$ws.Range("D21").Value = 104000000.0
$ws.Range("E21").Value = 62000000.0
$ws.Range("F21").Value = 66000000.0
$ws.Range("G21").Value = 76000000.0

# Row 32: Shares (Common) - fill previously empty B32
$ws.Range("B32").Value = 124360000.0

# Row 34: Net Debt - fill previously empty B34
$ws.Range("B34").Value = 2291000000.0

# Row 35: Total Debt - fill previously empty B35
$ws.Range("B35").Value = 2783000000.0
